$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 2.512729
$ws.Range("H2").Value = 7.538187
$ws.Range("I2").Value = 0.02190726325199687
$ws.Range("J2").Value = 0.02190726325199687
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 24.44779533333333
$ws.Range("N2").Value = 73.343386
$ws.Range("O2").Value = 0.1890645123346783
$ws.Range("P2").Value = 0.1890645123346783
$ws.Range("Q2").Value = 61.43068432013133
$ws.Range("R2").Value = 552.8761588811819
$ws.Range("S2").Value = 0.004141886043326206
$ws.Range("T2").Value = 0.004141886043326208

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 2.512729
$ws.Range("H3").Value = 7.538187
$ws.Range("I3").Value = 0.02190726325199687
$ws.Range("J3").Value = 0.02190726325199687
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 18.11074233333333
$ws.Range("N3").Value = 54.332227
$ws.Range("O3").Value = 0.1400575643155068
$ws.Range("P3").Value = 0.1400575643155068
$ws.Range("Q3").Value = 45.50738747249433
$ws.Range("R3").Value = 409.566487252449
$ws.Range("S3").Value = 0.003068277931893291
$ws.Range("T3").Value = 0.003068277931893291

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 2.512729
$ws.Range("H4").Value = 7.538187
$ws.Range("I4").Value = 0.02190726325199687
$ws.Range("J4").Value = 0.02190726325199687
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 80.684877
$ws.Range("N4").Value = 242.054631
$ws.Range("O4").Value = 0.6239682030546764
$ws.Range("P4").Value = 0.6239682030546765
$ws.Range("Q4").Value = 202.739230299333
$ws.Range("R4").Value = 1824.653072693997
$ws.Range("S4").Value = 0.01366943568519423
$ws.Range("T4").Value = 0.01366943568519424

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 2.512729
$ws.Range("H5").Value = 7.538187
$ws.Range("I5").Value = 0.02190726325199687
$ws.Range("J5").Value = 0.02190726325199687
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 6.065862
$ws.Range("N5").Value = 18.197586
$ws.Range("O5").Value = 0.0469097202951384
$ws.Range("P5").Value = 0.04690972029513841
$ws.Range("Q5").Value = 15.241867357398
$ws.Range("R5").Value = 137.176806216582
$ws.Range("S5").Value = 0.001027663591583137
$ws.Range("T5").Value = 0.001027663591583138

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 92.89399466666667
$ws.Range("H6").Value = 278.681984
$ws.Range("I6").Value = 0.8098976036382196
$ws.Range("J6").Value = 0.8098976036382197
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 24.44779533333333
$ws.Range("N6").Value = 73.343386
$ws.Range("O6").Value = 0.1890645123346783
$ws.Range("P6").Value = 0.1890645123346783
$ws.Range("Q6").Value = 2271.053369306425
$ws.Range("R6").Value = 20439.48032375782
$ws.Range("S6").Value = 0.1531228954728846
$ws.Range("T6").Value = 0.1531228954728846

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 92.89399466666667
$ws.Range("H7").Value = 278.681984
$ws.Range("I7").Value = 0.8098976036382196
$ws.Range("J7").Value = 0.8098976036382197
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 18.11074233333333
$ws.Range("N7").Value = 54.332227
$ws.Range("O7").Value = 0.1400575643155068
$ws.Range("P7").Value = 0.1400575643155068
$ws.Range("Q7").Value = 1682.379201722041
$ws.Range("R7").Value = 15141.41281549837
$ws.Range("S7").Value = 0.1134322857105348
$ws.Range("T7").Value = 0.1134322857105348

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 92.89399466666667
$ws.Range("H8").Value = 278.681984
$ws.Range("I8").Value = 0.8098976036382196
$ws.Range("J8").Value = 0.8098976036382197
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 80.684877
$ws.Range("N8").Value = 242.054631
$ws.Range("O8").Value = 0.6239682030546764
$ws.Range("P8").Value = 0.6239682030546765
$ws.Range("Q8").Value = 7495.140533718656
$ws.Range("R8").Value = 67456.2648034679
$ws.Range("S8").Value = 0.5053503524004285
$ws.Range("T8").Value = 0.5053503524004286

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 92.89399466666667
$ws.Range("H9").Value = 278.681984
$ws.Range("I9").Value = 0.8098976036382196
$ws.Range("J9").Value = 0.8098976036382197
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 6.065862
$ws.Range("N9").Value = 18.197586
$ws.Range("O9").Value = 0.0469097202951384
$ws.Range("P9").Value = 0.04690972029513841
$ws.Range("Q9").Value = 563.4821522767361
$ws.Range("R9").Value = 5071.339370490624
$ws.Range("S9").Value = 0.03799207005437175
$ws.Range("T9").Value = 0.03799207005437176

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 1.610639333333333
$ws.Range("H10").Value = 4.831918
$ws.Range("I10").Value = 0.0140423817607685
$ws.Range("J10").Value = 0.0140423817607685
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 24.44779533333333
$ws.Range("N10").Value = 73.343386
$ws.Range("O10").Value = 0.1890645123346783
$ws.Range("P10").Value = 0.1890645123346783
$ws.Range("Q10").Value = 39.37658077714978
$ws.Range("R10").Value = 354.389226994348
$ws.Range("S10").Value = 0.002654916059617078
$ws.Range("T10").Value = 0.002654916059617078

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 1.610639333333333
$ws.Range("H11").Value = 4.831918
$ws.Range("I11").Value = 0.0140423817607685
$ws.Range("J11").Value = 0.0140423817607685
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 18.11074233333333
$ws.Range("N11").Value = 54.332227
$ws.Range("O11").Value = 0.1400575643155068
$ws.Range("P11").Value = 0.1400575643155068
$ws.Range("Q11").Value = 29.16987395793178
$ws.Range("R11").Value = 262.528865621386
$ws.Range("S11").Value = 0.001966741786601734
$ws.Range("T11").Value = 0.001966741786601734

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 1.610639333333333
$ws.Range("H12").Value = 4.831918
$ws.Range("I12").Value = 0.0140423817607685
$ws.Range("J12").Value = 0.0140423817607685
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 80.684877
$ws.Range("N12").Value = 242.054631
$ws.Range("O12").Value = 0.6239682030546764
$ws.Range("P12").Value = 0.6239682030546765
$ws.Range("Q12").Value = 129.954236501362
$ws.Range("R12").Value = 1169.588128512258
$ws.Range("S12").Value = 0.008761999713874483
$ws.Range("T12").Value = 0.008761999713874487

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 1.610639333333333
$ws.Range("H13").Value = 4.831918
$ws.Range("I13").Value = 0.0140423817607685
$ws.Range("J13").Value = 0.0140423817607685
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 6.065862
$ws.Range("N13").Value = 18.197586
$ws.Range("O13").Value = 0.0469097202951384
$ws.Range("P13").Value = 0.04690972029513841
$ws.Range("Q13").Value = 9.769915927772001
$ws.Range("R13").Value = 87.92924334994801
$ws.Range("S13").Value = 0.0006587242006752034
$ws.Range("T13").Value = 0.0006587242006752036

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 17.68108066666667
$ws.Range("H14").Value = 53.04324200000001
$ws.Range("I14").Value = 0.154152751349015
$ws.Range("J14").Value = 0.154152751349015
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 24.44779533333333
$ws.Range("N14").Value = 73.343386
$ws.Range("O14").Value = 0.1890645123346783
$ws.Range("P14").Value = 0.1890645123346783
$ws.Range("Q14").Value = 432.2634414108236
$ws.Range("R14").Value = 3890.370972697412
$ws.Range("S14").Value = 0.02914481475885044
$ws.Range("T14").Value = 0.02914481475885045

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 17.68108066666667
$ws.Range("H15").Value = 53.04324200000001
$ws.Range("I15").Value = 0.154152751349015
$ws.Range("J15").Value = 0.154152751349015
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 18.11074233333333
$ws.Range("N15").Value = 54.332227
$ws.Range("O15").Value = 0.1400575643155068
$ws.Range("P15").Value = 0.1400575643155068
$ws.Range("Q15").Value = 320.2174961288816
$ws.Range("R15").Value = 2881.957465159935
$ws.Range("S15").Value = 0.021590258886477
$ws.Range("T15").Value = 0.021590258886477

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 17.68108066666667
$ws.Range("H16").Value = 53.04324200000001
$ws.Range("I16").Value = 0.154152751349015
$ws.Range("J16").Value = 0.154152751349015
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 80.684877
$ws.Range("N16").Value = 242.054631
$ws.Range("O16").Value = 0.6239682030546764
$ws.Range("P16").Value = 0.6239682030546765
$ws.Range("Q16").Value = 1426.595818817078
$ws.Range("R16").Value = 12839.3623693537
$ws.Range("S16").Value = 0.09618641525517922
$ws.Range("T16").Value = 0.09618641525517924

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 17.68108066666667
$ws.Range("H17").Value = 53.04324200000001
$ws.Range("I17").Value = 0.154152751349015
$ws.Range("J17").Value = 0.154152751349015
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 6.065862
$ws.Range("N17").Value = 18.197586
$ws.Range("O17").Value = 0.0469097202951384
$ws.Range("P17").Value = 0.04690972029513841
$ws.Range("Q17").Value = 107.250995334868
$ws.Range("R17").Value = 965.2589580138122
$ws.Range("S17").Value = 0.007231262448508312
$ws.Range("T17").Value = 0.007231262448508313
